$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11-16 down to 12-17
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with the latest week's data for this product/quality
$ws.Range("A11").Value = 3
$ws.Range("B11").Value = 'Femacal de La Calera'
$ws.Range("C11").Value = 'Coquimbo'
$ws.Range("D11").Value = 44992
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 100112043
$ws.Range("G11").Value = 'Pepino dulce'
$ws.Range("H11").Value = 'Cultivar IV Región'
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 56
$ws.Range("K11").Value = 13000
$ws.Range("L11").Value = 13000
$ws.Range("M11").Value = 13000
$ws.Range("N11").Value = '$/bandeja 18 kilos'
$ws.Range("O11").Value = 'Provincia de Limarí'
$ws.Range("P11").Value = 722
$ws.Range("Q11").Value = 18
$ws.Range("R11").Value = 'Hortaliza'
